$d = $word.ActiveDocument

# Locate the paragraphs that bracket the block being collapsed:
#  - the first "Section Headnote" paragraph, "What is a corporation?"
#  - the "Section Title" paragraph that follows it, "Section Two"
# Everything between them (the Resource/Case/Section paragraphs for
# "1.1 Case of the District Number 1" through "Section Two") gets folded
# into the first paragraph's text, dropping the second ResourceHeadnote's
# text ("This is an annotatable resource in the casebook."). The trailing
# "Section Headnote" paragraph ("This is the second chapter of the
# casebook.") is left untouched.
$startParaIndex = -1
$endParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $txt = $p.Range.Text.TrimEnd([char]13)
    if ($startParaIndex -lt 0 -and $p.Style.NameLocal -eq "Section Headnote" -and $txt -eq "What is a corporation?") {
        $startParaIndex = $i
    }
    if ($startParaIndex -gt 0 -and $endParaIndex -lt 0 -and $p.Style.NameLocal -eq "Section Title" -and $txt -eq "Section Two") {
        $endParaIndex = $i
    }
}

$startPara = $d.Paragraphs.Item($startParaIndex)
$endPara = $d.Paragraphs.Item($endParaIndex)

# Replace the whole multi-paragraph range (including paragraph marks) with
# a single run of plain text -- this both merges the paragraphs into one
# (keeping $startParaIndex's paragraph/style) and collapses all the runs
# into a single run, dropping the ResourceHeadnote text that isn't carried
# over.
$full = $d.Range($startPara.Range.Start, $endPara.Range.End)
$newText = "1.1Case of the District Number 1This is the body of case 1." + `
    "1.2Case of the District Number 2highlighted: content to highlight; " + `
    "elided: content to elide; replaced: content to replace; commented: " + `
    "content to comment; highlighted2: second highlight content;`n2Section Two"
$full.Text = $newText

# The assignment above leaves the now-empty paragraphs that used to hold
# "1.1", "Case of the District Number 1", etc. (they kept their paragraph
# marks). Remove that now-empty range to actually delete those paragraphs.
$mergedPara = $d.Paragraphs.Item($startParaIndex)
$nextPara = $d.Paragraphs.Item($endParaIndex)
$tail = $d.Range($mergedPara.Range.End, $nextPara.Range.End)
$tail.Delete()
